# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) for the fe5756af... row (row 3)
# on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-18 02:56:38"
$wsZhCn.Range("H3").Value = "2016-03-18 02:57:19"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-18 02:56:46"
$wsDeDe.Range("H3").Value = "2016-03-18 02:57:31"
